$d = $word.ActiveDocument

# 1. Remove the leading empty paragraph.
$d.Paragraphs.Item(1).Range.Delete()

# 2. Remove the trailing "+++END-FOR company+++" paragraph and the two
#    trailing empty paragraphs (everything after the LINK paragraph).
$startDel = $d.Paragraphs.Item(4).Range.Start
$endDel = $d.Paragraphs.Item(6).Range.End
$d.Range($startDel, $endDel).Delete()

# 3. "+++FOR company IN companies+++" -> "+++FOR link IN links+++"
$d.Content.Find.Execute("+++FOR company IN companies+++", $true, $false, $false, $false, $false, $true, 1, $false, "+++FOR link IN links+++", 2)

# 4. "+++INS `$company.name+++" -> "+++LINK ({  url: `$link.url, label: `$link.name })+++"
$d.Content.Find.Execute("+++INS `$company.name+++", $true, $false, $false, $false, $false, $true, 1, $false, "+++LINK ({  url: `$link.url, label: `$link.name })+++", 2)

# 5. "+++LINK ({ url: 'http://www.test.com', label: 'Link' })+++" -> "+++END-FOR link+++"
#    (the bookmark on this paragraph is preserved automatically)
$d.Content.Find.Execute("+++LINK ({ url: 'http://www.test.com', label: 'Link' })+++", $true, $false, $false, $false, $false, $true, 1, $false, "+++END-FOR link+++", 2)
